$d = $word.ActiveDocument

# 1. Update the letter date
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "September 21, 2025", 2)

# 2. Split the address line "979 Story Road, San Jose CA 95122" into two paragraphs:
#    "979 Story Road" and "San Jose, CA 95122"
$d.Content.Find.Execute("979 Story Road, San Jose CA 95122", $true, $false, $false, $false, $false,
                         $true, 1, $false, "979 Story Road^pSan Jose, CA 95122", 2)

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "San Jose, CA 95122`r") {
        $p.Range.Font.Name = "Arial"
        $p.Range.Font.NameAscii = "Arial"
        $p.Range.Font.NameBi = "Arial"
        $p.Range.Font.Size = 11
        $p.Range.Font.SizeBi = 11
        break
    }
}

# 3. Remove the empty NoSpacing paragraph that follows "Board of Directors"
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Vietnam Town Condominium Owners Association Board of Directors`r") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text -eq "`r") {
            $next.Range.Delete()
        }
        break
    }
}
